$d = $word.ActiveDocument

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>Chapter 12, Question 3</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:tab/>
        <w:t>This question was overridden, so we aren’t doing the question in the book. Instead, we are supposed to change the gravity from 9.8m/s</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
          <w:vertAlign w:val="superscript"/>
        </w:rPr>
        <w:t xml:space="preserve">2 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>to 3m/s</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
          <w:vertAlign w:val="superscript"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">. Additionally, we need to implement a slow down timer when the runner starts to fall and then resume the game’s regular speed when she stops falling. For the former, all we do is change the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>GRAVITY_FORCE</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> constant to 3.0. For the latter, we can use the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>setTimeRate</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> function inside of the </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>fall(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>stopFalling</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> functions. At the end of the </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>fall(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> function, we set the time rate to 0.25. At the end of the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>stopFalling</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Cascadia Code" w:hAnsi="Cascadia Code" w:cs="Cascadia Code"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> function, we set the time rate back to 1. This slows the game down while we’re falling and speeds it up once we land. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>This has the intended effect of slowing down the game when falling and returning the game back to normal speed after.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$end = $d.Content
$end.Collapse(0)
$end.InsertXML($xml) | Out-Null
